$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values per row (columns D:AJ) to reflect corrected IFRS figures
$ws.Range("D2").Value = 6700
$ws.Range("E2").Value = 337
$ws.Range("F2").Value = 337
$ws.Range("G2").Value = 465
$ws.Range("H2").Value = 371
$ws.Range("I2").Value = 238
$ws.Range("J2").Value = 133
$ws.Range("K2").Value = 8999
$ws.Range("L2").Value = 4600
$ws.Range("M2").Value = 4398
$ws.Range("N2").Value = 2811
$ws.Range("O2").Value = 1587
$ws.Range("P2").Value = 85
$ws.Range("Q2").Value = 529
$ws.Range("R2").Value = -241
$ws.Range("S2").Value = -82
$ws.Range("T2").Value = 191
$ws.Range("U2").Value = 338
$ws.Range("V2").Value = 2463
$ws.Range("W2").Value = 5.04
$ws.Range("X2").Value = 5.54
$ws.Range("Y2").Value = 8.68
$ws.Range("Z2").Value = 4.23
$ws.Range("AA2").Value = 104.59
$ws.Range("AB2").Value = 3378.67
$ws.Range("AC2").Value = 1765
$ws.Range("AD2").Value = 14.22
$ws.Range("AE2").Value = 26978
$ws.Range("AF2").Value = 0.93
$ws.Range("AH2").Value = 1.59
$ws.Range("AI2").Value = 17.53
$ws.Range("AJ2").Value = 13198611
$ws.Range("D3").Value = 6835
$ws.Range("E3").Value = 272
$ws.Range("F3").Value = 272
$ws.Range("G3").Value = 1119
$ws.Range("H3").Value = 882
$ws.Range("I3").Value = 644
$ws.Range("J3").Value = 238
$ws.Range("K3").Value = 9865
$ws.Range("L3").Value = 4651
$ws.Range("M3").Value = 5214
$ws.Range("N3").Value = 3422
$ws.Range("O3").Value = 1792
$ws.Range("P3").Value = 85
$ws.Range("Q3").Value = 358
$ws.Range("R3").Value = 698
$ws.Range("S3").Value = -261
$ws.Range("T3").Value = 249
$ws.Range("U3").Value = 109
$ws.Range("V3").Value = 2285
$ws.Range("W3").Value = 3.98
$ws.Range("X3").Value = 12.91
$ws.Range("Y3").Value = 20.67
$ws.Range("Z3").Value = 9.35
$ws.Range("AA3").Value = 89.20999999999999
$ws.Range("AB3").Value = 4077.28
$ws.Range("AC3").Value = 4780
$ws.Range("AD3").Value = 5.13
$ws.Range("AE3").Value = 32841
$ws.Range("AF3").Value = 0.75
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 1.63
$ws.Range("AI3").Value = 6.47
$ws.Range("AJ3").Value = 13198611
$ws.Range("D4").Value = 6901
$ws.Range("E4").Value = 227
$ws.Range("F4").Value = 227
$ws.Range("G4").Value = 349
$ws.Range("H4").Value = 209
$ws.Range("I4").Value = 109
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 9702
$ws.Range("L4").Value = 4367
$ws.Range("M4").Value = 5335
$ws.Range("N4").Value = 3492
$ws.Range("O4").Value = 1844
$ws.Range("P4").Value = 85
$ws.Range("Q4").Value = 545
$ws.Range("R4").Value = -621
$ws.Range("S4").Value = -371
$ws.Range("T4").Value = 257
$ws.Range("U4").Value = 288
$ws.Range("V4").Value = 2019
$ws.Range("W4").Value = 3.29
$ws.Range("X4").Value = 3.02
$ws.Range("Y4").Value = 3.16
$ws.Range("Z4").Value = 2.13
$ws.Range("AA4").Value = 81.84999999999999
$ws.Range("AB4").Value = 4151.56
$ws.Range("AC4").Value = 809
$ws.Range("AD4").Value = 19.28
$ws.Range("AE4").Value = 33509
$ws.Range("AF4").Value = 0.47
$ws.Range("AG4").Value = 450
$ws.Range("AH4").Value = 2.88
$ws.Range("AI4").Value = 43.01
$ws.Range("AJ4").Value = 13198611
$ws.Range("D5").Value = 7713
$ws.Range("E5").Value = 173
$ws.Range("F5").Value = 173
$ws.Range("G5").Value = 867
$ws.Range("H5").Value = 620
$ws.Range("I5").Value = 426
$ws.Range("J5").Value = 194
$ws.Range("K5").Value = 9886
$ws.Range("L5").Value = 4111
$ws.Range("M5").Value = 5775
$ws.Range("N5").Value = 3812
$ws.Range("O5").Value = 1963
$ws.Range("P5").Value = 85
$ws.Range("Q5").Value = 307
$ws.Range("R5").Value = 52
$ws.Range("S5").Value = -327
$ws.Range("T5").Value = 394
$ws.Range("U5").Value = -87
$ws.Range("V5").Value = 1791
$ws.Range("W5").Value = 2.24
$ws.Range("X5").Value = 8.039999999999999
$ws.Range("Y5").Value = 11.67
$ws.Range("Z5").Value = 6.33
$ws.Range("AA5").Value = 71.19
$ws.Range("AB5").Value = 4619.8
$ws.Range("AC5").Value = 3162
$ws.Range("AD5").Value = 5.09
$ws.Range("AE5").Value = 36584
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 450
$ws.Range("AH5").Value = 2.8
$ws.Range("AI5").Value = 11.01
$ws.Range("AJ5").Value = 13291151
$ws.Range("D6").Value = 8360
$ws.Range("E6").Value = 102
$ws.Range("F6").Value = 102
$ws.Range("G6").Value = 137
$ws.Range("H6").Value = 27
$ws.Range("I6").Value = -82
$ws.Range("K6").Value = 9871
$ws.Range("L6").Value = 4256
$ws.Range("M6").Value = 5615
$ws.Range("N6").Value = 3653
$ws.Range("P6").Value = 85
$ws.Range("Q6").Value = 74
$ws.Range("R6").Value = -176
$ws.Range("S6").Value = -52
$ws.Range("T6").Value = 175
$ws.Range("U6").Value = -101
$ws.Range("V6").Value = 1950
$ws.Range("W6").Value = 1.22
$ws.Range("X6").Value = 0.32
$ws.Range("Y6").Value = -2.19
$ws.Range("Z6").Value = 0.27
$ws.Range("AA6").Value = 75.81
$ws.Range("AB6").Value = 4470.41
$ws.Range("AC6").Value = -607
$ws.Range("AD6").Value = -20.17
$ws.Range("AE6").Value = 35053
$ws.Range("AF6").Value = 0.35
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 3.67
$ws.Range("AI6").Value = -57.31
$ws.Range("AJ6").Value = 13291151
$ws.Range("D7").Value = 9020
$ws.Range("E7").Value = 320
$ws.Range("G7").Value = 450
$ws.Range("H7").Value = 260
$ws.Range("I7").Value = 70
$ws.Range("K7").Value = 10360
$ws.Range("L7").Value = 4580
$ws.Range("M7").Value = 5770
$ws.Range("N7").Value = 3710
$ws.Range("P7").Value = 90
$ws.Range("Q7").Value = 580
$ws.Range("R7").Value = -520
$ws.Range("S7").Value = -130
$ws.Range("T7").Value = 220
$ws.Range("W7").Value = 3.55
$ws.Range("X7").Value = 2.88
$ws.Range("Y7").Value = 1.9
$ws.Range("Z7").Value = 2.57
$ws.Range("AA7").Value = 79.38
$ws.Range("AC7").Value = 519
$ws.Range("AD7").Value = 20.89
$ws.Range("AE7").Value = 35603
$ws.Range("AF7").Value = 0.3
$ws.Range("AG7").Value = 450
$ws.Range("AH7").Value = 4.15
$ws.Range("AI7").Value = 85.44
$ws.Range("D8").Value = 9460
$ws.Range("E8").Value = 390
$ws.Range("G8").Value = 520
$ws.Range("H8").Value = 300
$ws.Range("I8").Value = 110
$ws.Range("K8").Value = 10590
$ws.Range("L8").Value = 4560
$ws.Range("M8").Value = 6030
$ws.Range("N8").Value = 3770
$ws.Range("P8").Value = 90
$ws.Range("Q8").Value = 390
$ws.Range("R8").Value = -240
$ws.Range("S8").Value = -50
$ws.Range("T8").Value = 200
$ws.Range("W8").Value = 4.12
$ws.Range("X8").Value = 3.17
$ws.Range("Y8").Value = 2.94
$ws.Range("Z8").Value = 2.86
$ws.Range("AA8").Value = 75.62
$ws.Range("AC8").Value = 816
$ws.Range("AD8").Value = 13.29
$ws.Range("AE8").Value = 36179
$ws.Range("AF8").Value = 0.3
$ws.Range("AG8").Value = 450
$ws.Range("AH8").Value = 4.15
$ws.Range("AI8").Value = 54.37
$ws.Range("D9").Value = 10090
$ws.Range("E9").Value = 440
$ws.Range("G9").Value = 600
$ws.Range("H9").Value = 350
$ws.Range("I9").Value = 120
$ws.Range("K9").Value = 10970
$ws.Range("L9").Value = 4640
$ws.Range("M9").Value = 6330
$ws.Range("N9").Value = 3850
$ws.Range("P9").Value = 90
$ws.Range("Q9").Value = 400
$ws.Range("R9").Value = -310
$ws.Range("S9").Value = 30
$ws.Range("T9").Value = 200
$ws.Range("W9").Value = 4.36
$ws.Range("X9").Value = 3.47
$ws.Range("Y9").Value = 3.15
$ws.Range("Z9").Value = 3.25
$ws.Range("AA9").Value = 73.3
$ws.Range("AC9").Value = 890
$ws.Range("AD9").Value = 12.18
$ws.Range("AE9").Value = 36946
$ws.Range("AF9").Value = 0.29
$ws.Range("AG9").Value = 450
$ws.Range("AH9").Value = 4.15
$ws.Range("AI9").Value = 49.84

# Clear cells that no longer contain data in the corrected sheet
$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()
